$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 110, pushing existing rows 110:156 down to 111:157
$ws.Rows.Item(110).Insert()

# Copy the static (unchanged) columns from the row that is now 111 (old row 110)
$staticCols = @("A","B","C","E","F","G","H","I","J","Q","R","T")
foreach ($col in $staticCols) {
    $src = $col + "111"
    $dst = $col + "110"
    $ws.Range($dst).Value = $ws.Range($src).Value()
}

# Set the new data values for row 110
$ws.Range("D110").Value = 44489
$ws.Range("K110").Value = "Red Blush"
$ws.Range("L110").Value = "Primera"
$ws.Range("M110").Value = 30
$ws.Range("N110").Value = 11000
$ws.Range("O110").Value = 12000
$ws.Range("P110").Value = 11500
$ws.Range("S110").Value = 821
